$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.566.11'
$ws.Range("E2").Value = '  +2.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.474.37'
$ws.Range("E3").Value = '  +2.25%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.00'
$ws.Range("E5").Value = '  +2.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.47'
$ws.Range("E6").Value = '  +3.68%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.542'
$ws.Range("E8").Value = '  +1.91%  '
$ws.Range("E9").Value = '  +4.60%  '
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.364'
$ws.Range("E11").Value = '  +4.09%  '
$ws.Range("E12").Value = '  +2.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.33'
$ws.Range("E13").Value = '  +4.18%  '
$ws.Range("E14").Value = '  +6.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.386.87'
$ws.Range("E16").Value = '  +2.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.481.22'
$ws.Range("E17").Value = '  +2.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.57'
$ws.Range("E18").Value = '  +2.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.31'
$ws.Range("E19").Value = '  +7.00%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.25'
$ws.Range("E20").Value = '  +2.70%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '329.21'
$ws.Range("E21").Value = '  +1.77%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("E23").Value = '  +10.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.49'
$ws.Range("E24").Value = '  +1.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '631.78'
$ws.Range("E25").Value = '  +14.05%  '
$ws.Range("E26").Value = '  +13.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.76'
$ws.Range("E27").Value = '  -0.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.595.75'
$ws.Range("E28").Value = '  +2.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.53'
$ws.Range("E29").Value = '  +9.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.47'
$ws.Range("E30").Value = '  +3.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.997'
$ws.Range("E31").Value = '  -0.51%  '
$ws.Range("E32").Value = '  -1.81%  '
$ws.Range("E33").Value = '  +1.65%  '
$ws.Range("E34").Value = '  +9.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.57'
$ws.Range("E35").Value = '  +3.94%  '
$ws.Range("E36").Value = '  -0.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.387'
$ws.Range("E37").Value = '  +2.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.55'
$ws.Range("E38").Value = '  +2.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.99'
$ws.Range("E39").Value = '  +2.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.86'
$ws.Range("E40").Value = '  +2.17%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.72'
$ws.Range("E41").Value = '  +21.93%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '147.18'
$ws.Range("E42").Value = '  -3.80%  '
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '150.95'
$ws.Range("E44").Value = '  +2.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.78'
$ws.Range("E45").Value = '  +3.71%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.28'
$ws.Range("E46").Value = '  +7.37%  '
$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0552'
$ws.Range("E47").Value = '  +4.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.611'
$ws.Range("E48").Value = '  +3.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0241'
$ws.Range("E49").Value = '  +5.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0930'
$ws.Range("E50").Value = '  +1.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.750'
$ws.Range("E51").Value = '  +5.15%  '
